$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.962.19"
$ws.Range("E2").Value = "  -1.69%  "
$ws.Range("D3").Value = "3.298.38"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.61"
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.54"
$ws.Range("E6").Value = "  -4.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.629"
$ws.Range("E7").Value = "  +4.73%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.126"
$ws.Range("E9").Value = "  -2.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.68"
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.403"
$ws.Range("E11").Value = "  -1.97%  "
$ws.Range("D12").Value = "3.876.39"
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("E13").Value = "  -3.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.63"
$ws.Range("E14").Value = "  -3.09%  "
$ws.Range("D15").Value = "66.093.32"
$ws.Range("E15").Value = "  -1.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000164"
$ws.Range("E16").Value = "  -1.93%  "
$ws.Range("D17").Value = "3.308.75"
$ws.Range("E17").Value = "  +1.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "435.85"
$ws.Range("E18").Value = "  -1.86%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.61"
$ws.Range("E19").Value = "  -1.83%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.38"
$ws.Range("E20").Value = "  -1.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.47"
$ws.Range("E21").Value = "  -3.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.85"
$ws.Range("E22").Value = "  -2.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.516"
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("D25").Value = "3.461.78"
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000114"
$ws.Range("E26").Value = "  -4.01%  "
$ws.Range("E27").Value = "  +4.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.01"
$ws.Range("E28").Value = "  -1.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.94"
$ws.Range("E30").Value = "  -1.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.54"
$ws.Range("E31").Value = "  -1.28%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.17"
$ws.Range("E33").Value = "  -3.54%  "
$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.69"
$ws.Range("E34").Value = "  -1.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.20"
$ws.Range("E35").Value = "  -3.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.46"
$ws.Range("E36").Value = "  -4.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.03"
$ws.Range("E37").Value = "  -2.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "27.23"
$ws.Range("E38").Value = "  -0.91%  "
$ws.Range("E39").Value = "  -2.68%  "
$ws.Range("D40").Value = "2.790.66"
$ws.Range("E40").Value = "  +2.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.784"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.37"
$ws.Range("E42").Value = "  -2.40%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.13"
$ws.Range("E43").Value = "  -2.44%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.55"
$ws.Range("E44").Value = "  +0.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0660"
$ws.Range("E45").Value = "  -1.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "323.66"
$ws.Range("E46").Value = "  -1.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.30"
$ws.Range("E47").Value = "  -4.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.60"
$ws.Range("E48").Value = "  -4.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0271"
$ws.Range("E49").Value = "  -1.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.102"
$ws.Range("E50").Value = "  +2.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.10"
$ws.Range("E51").Value = "  -2.04%  "
